$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared timestamp used by the first two test cases (rows 2-3) is
# refreshed to the new execution run.
$ws.Range("D2").Value = "13/05/2025 11:04:40 AM"
$ws.Range("D3").Value = "13/05/2025 11:04:40 AM"

# Row 4 ("Navigate to Profile Page") now records the execution time of
# the new run (one second later than the shared timestamp above).
$ws.Range("D4").Value = "13/05/2025 11:04:41 AM"

# A new test case was added: "Verify profile page loads".
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = " Verify profile page loads"
$ws.Range("C5").Value = "PASSED"
$ws.Range("D5").Value = "13/05/2025 11:04:41 AM"
$ws.Range("E5").Value = "Test executed successfully."
